# fct_portfolio_pnl.xlsx — refresh of exported P&L metrics.
# The upstream recompute (new customer loan-count logic + CAC-per-loan
# amortization feeding acquisition_cost) shifts fully_paid_loans (Z) for
# every vintage/segment row, which in turn nudges every downstream ratio
# that derives from it (margins, per-loan economics, LTV:CAC) by the
# corresponding floating-point amount. Apply the refreshed cell values
# row by row, matching the regenerated export exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 — 2025-01-01 / High Risk
$ws.Range("I2").Value  = 258500.46
$ws.Range("M2").Value  = 258500.46
$ws.Range("N2").Value  = 258120.0084129708
$ws.Range("O2").Value  = 223414.3624129708
$ws.Range("R2").Value  = 177.1585507295613
$ws.Range("T2").Value  = 0.001471763675117404
$ws.Range("V2").Value  = 0.9985282363248825
$ws.Range("W2").Value  = 0.8642706570540372
$ws.Range("Z2").Value  = 1140
$ws.Range("AB2").Value = 7.437406824612077

# Row 3 — 2025-01-01 / Low Risk
$ws.Range("F3").Value = 7740.429999999999
$ws.Range("Z3").Value = 103

# Row 4 — 2025-01-01 / Medium Risk
$ws.Range("F4").Value = 84780.06999999999
$ws.Range("T4").Value = 0.0002695395723226645
$ws.Range("W4").Value = 0.7624469053905757
$ws.Range("Z4").Value = 547

# Row 6 — 2025-01-01 / Unknown
$ws.Range("Z6").Value = 8

# Row 7 — 2025-02-01 / High Risk
$ws.Range("E7").Value  = 2545370.780000001
$ws.Range("F7").Value  = 340922.6199999999
$ws.Range("I7").Value  = 496462.91
$ws.Range("K7").Value  = 508.2833202709305
$ws.Range("L7").Value  = 52107.251
$ws.Range("M7").Value  = 496462.91
$ws.Range("N7").Value  = 495954.6266797292
$ws.Range("O7").Value  = 443847.3756797291
$ws.Range("R7").Value  = 209.4402984289397
$ws.Range("V7").Value  = 0.9989761907485277
$ws.Range("W7").Value  = 0.8940192041329513
$ws.Range("Z7").Value  = 2003
$ws.Range("AB7").Value = 9.517958003190941

# Row 8 — 2025-02-01 / Low Risk
$ws.Range("L8").Value  = 9008.424500000001
$ws.Range("W8").Value  = 0.4844229878821072
$ws.Range("Z8").Value  = 224
$ws.Range("AB8").Value = 1.939574450560139

# Row 9 — 2025-02-01 / Medium Risk
$ws.Range("L9").Value = 41511.6055
$ws.Range("N9").Value = 230998.0738730376
$ws.Range("R9").Value = 213.8871054379978
$ws.Range("W9").Value = 0.8201865145260109
$ws.Range("Z9").Value = 1062

# Row 10 — 2025-02-01 / Other
$ws.Range("I10").Value = 813.3199999999999
$ws.Range("M10").Value = 813.3199999999999
$ws.Range("N10").Value = 813.3199999999999
$ws.Range("O10").Value = 813.3199999999999
$ws.Range("Z10").Value = 4

# Row 11 — 2025-02-01 / Unknown
$ws.Range("E11").Value = 31374.71999999999
$ws.Range("I11").Value = 4404.940000000001
$ws.Range("M11").Value = 4404.940000000001
$ws.Range("N11").Value = 4404.940000000001
$ws.Range("Z11").Value = 18

# Row 12 — 2025-03-01 / High Risk
$ws.Range("F12").Value = 301217.22
$ws.Range("G12").Value = 76912.54999999999
$ws.Range("N12").Value = 426450.1746825941
$ws.Range("P12").Value = 198.6544228993536
$ws.Range("S12").Value = 0.008736560201845529
$ws.Range("U12").Value = 0.9912634397981546
$ws.Range("Z12").Value = 1936

# Row 13 — 2025-03-01 / Low Risk
$ws.Range("L13").Value = 4709.800499999999
$ws.Range("Z13").Value = 240

# Row 14 — 2025-03-01 / Medium Risk
$ws.Range("F14").Value  = 168572.54
$ws.Range("K14").Value  = 30.43612696233116
$ws.Range("O14").Value  = 181912.4013730376
$ws.Range("V14").Value  = 0.99985397431846
$ws.Range("Z14").Value  = 1111
$ws.Range("AB14").Value = 7.867956357774901

# Row 15 — 2025-03-01 / Unknown
$ws.Range("Z15").Value  = 22
$ws.Range("AB15").Value = 7.535795727022715
